$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below mirror the upstream "cryptos" refresh run: updated Price (D) /
# Volume(1h) (E) readings, plus a swap of the KickToken / BKEXToken rows (41-42)
# including their Coin name (B) and Link (C).
#
# D/E hold numeric- and percent-looking text (e.g. "256.73", "0.43%") that must
# stay plain strings, exactly as they were before the edit -- so each of those
# cells is flipped to Text format ("@") before the assignment (stopping Excel
# from coercing "256.73" -> a number or "0.43%" -> 0.0043) and back to General
# afterwards to keep the original formatting/no-quote-prefix appearance.

$updates = [ordered]@{
    "D2" = "256.73"
    "E2" = "0.43%"
    "D3" = "27.06"
    "E3" = "-4.01%"
    "D4" = "4.626"
    "E4" = "-11.08%"
    "E5" = "0.57%"
    "D6" = "6.633"
    "E6" = "-0.80%"
    "D7" = "0.8659"
    "E7" = "-0.46%"
    "D8" = "0.9368"
    "E8" = "-3.14%"
    "D9" = "0.1402"
    "E9" = "-0.42%"
    "D10" = "0.03860"
    "E10" = "10.60%"
    "D11" = "0.07076"
    "E11" = "-0.63%"
    "D12" = "0.03198"
    "E12" = "0.60%"
    "D13" = "0.09248"
    "E13" = "0.40%"
    "D14" = "0.001540"
    "E14" = "-0.01%"
    "D15" = "0.0006043"
    "E15" = "-0.57%"
    "D16" = "0.006007"
    "E16" = "0.67%"
    "D17" = "3.514"
    "E17" = "0.52%"
    "D18" = "3.189"
    "D19" = "2.212"
    "E19" = "-0.58%"
    "D20" = "0.3099"
    "E20" = "-2.30%"
    "D21" = "0.1275"
    "E21" = "-2.57%"
    "D22" = "3.854"
    "E22" = "9.34%"
    "D23" = "0.04225"
    "E23" = "0.84%"
    "D24" = "0.001218"
    "E24" = "-0.58%"
    "D25" = "0.004283"
    "E25" = "-5.93%"
    "E26" = "-0.06%"
    "D27" = "0.0001935"
    "D40" = "0.03831"
    "E40" = "0.36%"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "D41" = "0.1101"
    "E41" = "-0.05%"
    "B42" = "KickToken"
    "C42" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "D42" = "0.003953"
    "E42" = "-29.53%"
    "D43" = "0.002312"
    "E43" = "-1.34%"
    "D44" = "0.01136"
    "E44" = "16.98%"
    "D45" = "0.00005447"
    "E45" = "1.45%"
    "E46" = "-0.09%"
    "D47" = "0.07771"
    "E47" = "-18.17%"
    "D48" = "0.002276"
    "E48" = "6.87%"
    "D49" = "0.00002098"
    "E49" = "-0.09%"
    "D50" = "0.0001998"
    "E50" = "-0.09%"
}

foreach ($addr in $updates.Keys) {
    $col = $addr -replace '[0-9]+$', ''
    $range = $ws.Range($addr)
    if ($col -eq "D" -or $col -eq "E") {
        # Numeric-looking text -> force Text format so it is not reinterpreted
        # as a number/percentage, then restore General (matches the source file).
        $range.NumberFormat = "@"
        $range.Value = $updates[$addr]
        $range.NumberFormat = "General"
    } else {
        $range.Value = $updates[$addr]
    }
}

Write-Host "Applied $($updates.Count) cell updates"
